$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the cell value: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select E8 so it becomes the active cell (matches saved selection in diff)
$ws.Activate()
$ws.Range("E8").Select()
